# Apply the changes described by the diff:
#   Trade #124 closed at 2026-02-16 21:46:27 - leadlag UP +0.000%
#
# Helper function: write a text value to a cell while avoiding Excel's
# automatic reinterpretation of strings that look like numbers/percents/dates
# (which would otherwise attach a numeric/date style to the cell). We force
# the cell to Text format, assign the literal string, then restore the
# cell's style to Normal so no stray style index is left behind.
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: update OVERALL and leadlag aggregate rows
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("C2").Value = 106
Set-TextValue $summary.Range("D2") "72.6%"
Set-TextValue $summary.Range("E2") "+34.6537%"
Set-TextValue $summary.Range("F2") "+0.3269%"

$summary.Range("C3").Value = 97
Set-TextValue $summary.Range("D3") "55.7%"
Set-TextValue $summary.Range("E3") "+20.8431%"
Set-TextValue $summary.Range("F3") "+0.2149%"

# ---------------------------------------------------------------------
# leadlag sheet: trades #104-106 (rows 80-82) move from OPEN to CLOSED,
# and a brand-new trade #124 (row 99) is appended as OPEN.
# ---------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

# Row 80 -> trade #104, now CLOSED
$leadlag.Range("G80").Value = 68885.167933
$leadlag.Range("H80").Value = "CLOSED"
$leadlag.Range("I80").Value = 0.5931
$leadlag.Range("J80").Value = 5.93
Set-TextValue $leadlag.Range("M80") "time_exit_5min"
$leadlag.Range("N80").Value = 5

# Row 81 -> trade #105, now CLOSED
$leadlag.Range("G81").Value = 68505.45422499999
$leadlag.Range("H81").Value = "CLOSED"
$leadlag.Range("I81").Value = -0.1606
$leadlag.Range("J81").Value = -1.61
Set-TextValue $leadlag.Range("M81") "time_exit_5min"
$leadlag.Range("N81").Value = 5

# Row 82 -> trade #106, now CLOSED
$leadlag.Range("G82").Value = 67574.21526
$leadlag.Range("H82").Value = "CLOSED"
$leadlag.Range("I82").Value = 1.1799
$leadlag.Range("J82").Value = 11.8
Set-TextValue $leadlag.Range("M82") "time_exit_5min"
$leadlag.Range("N82").Value = 5

# New row 99 -> trade #124, newly opened trade
$leadlag.Range("A99").Value = 124
Set-TextValue $leadlag.Range("B99") "2026-02-16"
Set-TextValue $leadlag.Range("C99") "21:46:27"
Set-TextValue $leadlag.Range("D99") "leadlag"
Set-TextValue $leadlag.Range("E99") "UP"
$leadlag.Range("F99").Value = 68473.485
$leadlag.Range("H99").Value = "OPEN"
$leadlag.Range("I99").Value = 0
$leadlag.Range("J99").Value = 0
$leadlag.Range("K99").Value = 0.6772
Set-TextValue $leadlag.Range("L99") "Binance leading with 0.068% move"

$leadlag.Range("N99").Value = 0

# ---------------------------------------------------------------------
# All Trades sheet: append the three now-CLOSED trades as new rows
# 105-107 (duplicated from the leadlag sheet, as this sheet aggregates
# across all strategies).
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Row 105 -> trade #104
$allTrades.Range("A105").Value = 104
Set-TextValue $allTrades.Range("B105") "2026-02-16"
Set-TextValue $allTrades.Range("C105") "21:40:50"
Set-TextValue $allTrades.Range("D105") "leadlag"
Set-TextValue $allTrades.Range("E105") "UP"
$allTrades.Range("F105").Value = 68479.035
$allTrades.Range("G105").Value = 68885.167933
Set-TextValue $allTrades.Range("H105") "CLOSED"
$allTrades.Range("I105").Value = 0.5931
$allTrades.Range("J105").Value = 5.93
$allTrades.Range("K105").Value = 0.75
Set-TextValue $allTrades.Range("L105") "Binance leading with 0.104% move"
Set-TextValue $allTrades.Range("M105") "time_exit_5min"
$allTrades.Range("N105").Value = 5

# Row 106 -> trade #105
$allTrades.Range("A106").Value = 105
Set-TextValue $allTrades.Range("B106") "2026-02-16"
Set-TextValue $allTrades.Range("C106") "21:41:02"
Set-TextValue $allTrades.Range("D106") "leadlag"
Set-TextValue $allTrades.Range("E106") "DOWN"
$allTrades.Range("F106").Value = 68395.58
$allTrades.Range("G106").Value = 68505.45422499999
Set-TextValue $allTrades.Range("H106") "CLOSED"
$allTrades.Range("I106").Value = -0.1606
$allTrades.Range("J106").Value = -1.61
$allTrades.Range("K106").Value = 0.75
Set-TextValue $allTrades.Range("L106") "Binance leading with -0.126% move"
Set-TextValue $allTrades.Range("M106") "time_exit_5min"
$allTrades.Range("N106").Value = 5

# Row 107 -> trade #106
$allTrades.Range("A107").Value = 106
Set-TextValue $allTrades.Range("B107") "2026-02-16"
Set-TextValue $allTrades.Range("C107") "21:41:08"
Set-TextValue $allTrades.Range("D107") "leadlag"
Set-TextValue $allTrades.Range("E107") "DOWN"
$allTrades.Range("F107").Value = 68381.06
$allTrades.Range("G107").Value = 67574.21526
Set-TextValue $allTrades.Range("H107") "CLOSED"
$allTrades.Range("I107").Value = 1.1799
$allTrades.Range("J107").Value = 11.8
$allTrades.Range("K107").Value = 0.75
Set-TextValue $allTrades.Range("L107") "Coinbase leading with -0.121% move"
Set-TextValue $allTrades.Range("M107") "time_exit_5min"
$allTrades.Range("N107").Value = 5

# ---------------------------------------------------------------------
# Comparison sheet: update leadlag aggregate row
# ---------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

$comparison.Range("B2").Value = 97
Set-TextValue $comparison.Range("C2") "55.7%"
Set-TextValue $comparison.Range("D2") "3.64"
Set-TextValue $comparison.Range("E2") "+0.5324%"
Set-TextValue $comparison.Range("F2") "-0.2928%"
Set-TextValue $comparison.Range("G2") "1.82"

Write-Host "Edit applied successfully"
